$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Row 20: "New Look and Feel" module is now checked off as done.
# Flip both the linked cell and the form-control's own state.
$ws.Range("J20").Value = $true
$chk = $ws.Shapes.Item("Check Box 32")
$chk.ControlFormat.Value = 1

# Row 37: task description changed from "Come up with your own" to "lock and key",
# and its point value increased from 1 to 2.
$ws.Range("C37").Value = "lock and key"
$ws.Range("D37").Value = 2

# Update the saved view: scrolled down a bit and a new active selection.
$excel.ActiveWindow.TopLeftCell = $ws.Range("A10")
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("O27").Select()
